# Scheduled market-data refresh: update currentAveragePrice / LevePrice / LeveProfit
# columns (H-N) across the ALC, ARM, BSM, CRP, CUL, GSM, LTW and WVR sheets.
# Columns: H=currentAveragePrice, I=currentAveragePriceNQ, J=currentAveragePriceHQ,
#          K=LevePriceNQ, L=LevePriceHQ, M=LeveProfitNQ, N=LeveProfitHQ
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC row 15
$ws.Cells.Item(15, 8).Value = 4463.1377
$ws.Cells.Item(15, 9).Value = 4463.1377
$ws.Cells.Item(15, 11).Value = 13389.4131
$ws.Cells.Item(15, 13).Value = -13220.4131

# ALC row 41
$ws.Cells.Item(41, 8).Value = 400.2
$ws.Cells.Item(41, 9).Value = 400
$ws.Cells.Item(41, 10).Value = 400.2857
$ws.Cells.Item(41, 11).Value = 400
$ws.Cells.Item(41, 12).Value = 400.2857
$ws.Cells.Item(41, 13).Value = 40
$ws.Cells.Item(41, 14).Value = -1280.2857

# ALC row 95
$ws.Cells.Item(95, 8).Value = 23972.25
$ws.Cells.Item(95, 10).Value = 23972.25
$ws.Cells.Item(95, 12).Value = 23972.25
$ws.Cells.Item(95, 14).Value = -29464.25

# ALC row 107
$ws.Cells.Item(107, 8).Value = 431.9
$ws.Cells.Item(107, 9).Value = 458.1111
$ws.Cells.Item(107, 10).Value = 196
$ws.Cells.Item(107, 11).Value = 458.1111
$ws.Cells.Item(107, 12).Value = 196
$ws.Cells.Item(107, 13).Value = 1461.8889
$ws.Cells.Item(107, 14).Value = -4036

# ALC row 138
$ws.Cells.Item(138, 8).Value = 3500192.2
$ws.Cells.Item(138, 10).Value = 3975812.5
$ws.Cells.Item(138, 12).Value = 11927437.5
$ws.Cells.Item(138, 14).Value = -11937717.5

# ALC row 139
$ws.Cells.Item(139, 8).Value = 34945
$ws.Cells.Item(139, 10).Value = 34945
$ws.Cells.Item(139, 12).Value = 34945
$ws.Cells.Item(139, 14).Value = -45225

$ws = $wb.Worksheets.Item("ARM")
# ARM row 2
$ws.Cells.Item(2, 8).Value = 1546.8334
$ws.Cells.Item(2, 9).Value = 1614.7273
$ws.Cells.Item(2, 10).Value = 800
$ws.Cells.Item(2, 11).Value = 1614.7273
$ws.Cells.Item(2, 12).Value = 800
$ws.Cells.Item(2, 13).Value = -1501.7273
$ws.Cells.Item(2, 14).Value = -1026

# ARM row 32
$ws.Cells.Item(32, 8).Value = 5669841
$ws.Cells.Item(32, 9).Value = 6826050
$ws.Cells.Item(32, 11).Value = 6826050
$ws.Cells.Item(32, 13).Value = -6825763

# ARM row 44
$ws.Cells.Item(44, 8).Value = 0
$ws.Cells.Item(44, 10).Value = 0
$ws.Cells.Item(44, 12).Value = 0
$ws.Cells.Item(44, 14).ClearContents()  # N44: clear (no longer reported)

# ARM row 61
$ws.Cells.Item(61, 8).Value = 32324866
$ws.Cells.Item(61, 9).Value = 45502036
$ws.Cells.Item(61, 10).Value = 114000
$ws.Cells.Item(61, 11).Value = 45502036
$ws.Cells.Item(61, 12).Value = 114000
$ws.Cells.Item(61, 13).Value = -45501824
$ws.Cells.Item(61, 14).Value = -114424

# ARM row 116
$ws.Cells.Item(116, 8).Value = 1546.8334
$ws.Cells.Item(116, 9).Value = 1614.7273
$ws.Cells.Item(116, 10).Value = 800
$ws.Cells.Item(116, 11).Value = 1614.7273
$ws.Cells.Item(116, 12).Value = 800
$ws.Cells.Item(116, 13).Value = 679.2727
$ws.Cells.Item(116, 14).Value = -5388

# ARM row 136
$ws.Cells.Item(136, 8).Value = 32324866
$ws.Cells.Item(136, 9).Value = 45502036
$ws.Cells.Item(136, 10).Value = 114000
$ws.Cells.Item(136, 11).Value = 136506108
$ws.Cells.Item(136, 12).Value = 342000
$ws.Cells.Item(136, 13).Value = -136503558
$ws.Cells.Item(136, 14).Value = -347100

$ws = $wb.Worksheets.Item("BSM")
# BSM row 3
$ws.Cells.Item(3, 8).Value = 1546.8334
$ws.Cells.Item(3, 9).Value = 1614.7273
$ws.Cells.Item(3, 10).Value = 800
$ws.Cells.Item(3, 11).Value = 1614.7273
$ws.Cells.Item(3, 12).Value = 800
$ws.Cells.Item(3, 13).Value = -1500.7273
$ws.Cells.Item(3, 14).Value = -1028

# BSM row 114
$ws.Cells.Item(114, 8).Value = 36842
$ws.Cells.Item(114, 10).Value = 36842
$ws.Cells.Item(114, 12).Value = 36842
$ws.Cells.Item(114, 14).Value = -45520

# BSM row 134
$ws.Cells.Item(134, 8).Value = 4375.0347
$ws.Cells.Item(134, 9).Value = 3726.353
$ws.Cells.Item(134, 10).Value = 5294
$ws.Cells.Item(134, 11).Value = 11179.059
$ws.Cells.Item(134, 12).Value = 15882
$ws.Cells.Item(134, 13).Value = -8644.059000000001
$ws.Cells.Item(134, 14).Value = -20952

$ws = $wb.Worksheets.Item("CRP")
# CRP row 16
$ws.Cells.Item(16, 8).Value = 1846.8334
$ws.Cells.Item(16, 9).Value = 1520.25
$ws.Cells.Item(16, 10).Value = 2500
$ws.Cells.Item(16, 11).Value = 1520.25
$ws.Cells.Item(16, 12).Value = 2500
$ws.Cells.Item(16, 13).Value = -1233.25
$ws.Cells.Item(16, 14).Value = -3074

# CRP row 113
$ws.Cells.Item(113, 8).Value = 1846.8334
$ws.Cells.Item(113, 9).Value = 1520.25
$ws.Cells.Item(113, 10).Value = 2500
$ws.Cells.Item(113, 11).Value = 1520.25
$ws.Cells.Item(113, 12).Value = 2500
$ws.Cells.Item(113, 13).Value = 649.75
$ws.Cells.Item(113, 14).Value = -6840

$ws = $wb.Worksheets.Item("CUL")
# CUL row 68
$ws.Cells.Item(68, 8).Value = 887.36664
$ws.Cells.Item(68, 10).Value = 952.9583
$ws.Cells.Item(68, 12).Value = 2858.8749
$ws.Cells.Item(68, 14).Value = -4480.8749

# CUL row 71
$ws.Cells.Item(71, 8).Value = 887.36664
$ws.Cells.Item(71, 10).Value = 952.9583
$ws.Cells.Item(71, 12).Value = 8576.6247
$ws.Cells.Item(71, 14).Value = -16688.6247

# CUL row 107
$ws.Cells.Item(107, 8).Value = 1046.3846
$ws.Cells.Item(107, 9).Value = 541.3125
$ws.Cells.Item(107, 10).Value = 1270.8611
$ws.Cells.Item(107, 11).Value = 1623.9375
$ws.Cells.Item(107, 12).Value = 3812.5833
$ws.Cells.Item(107, 13).Value = 296.0625
$ws.Cells.Item(107, 14).Value = -7652.5833

# CUL row 131
$ws.Cells.Item(131, 8).Value = 910.46155
$ws.Cells.Item(131, 9).Value = 571.3333
$ws.Cells.Item(131, 10).Value = 972.1212
$ws.Cells.Item(131, 11).Value = 1713.9999
$ws.Cells.Item(131, 12).Value = 2916.3636
$ws.Cells.Item(131, 13).Value = 3326.0001
$ws.Cells.Item(131, 14).Value = -12996.3636

$ws = $wb.Worksheets.Item("GSM")
# GSM row 29
$ws.Cells.Item(29, 8).Value = 25050
$ws.Cells.Item(29, 9).Value = 100
$ws.Cells.Item(29, 10).Value = 50000
$ws.Cells.Item(29, 11).Value = 100
$ws.Cells.Item(29, 12).Value = 50000
$ws.Cells.Item(29, 13).Value = 190
$ws.Cells.Item(29, 14).Value = -50580

# GSM row 107
$ws.Cells.Item(107, 8).Value = 4000
$ws.Cells.Item(107, 9).Value = 0
$ws.Cells.Item(107, 10).Value = 4000
$ws.Cells.Item(107, 11).Value = 0
$ws.Cells.Item(107, 12).Value = 4000
$ws.Cells.Item(107, 13).ClearContents()  # M107: clear (no longer reported)
$ws.Cells.Item(107, 14).Value = -7840

# GSM row 132
$ws.Cells.Item(132, 8).Value = 253624.62
$ws.Cells.Item(132, 9).Value = 252750
$ws.Cells.Item(132, 10).Value = 254499.25
$ws.Cells.Item(132, 11).Value = 758250
$ws.Cells.Item(132, 12).Value = 763497.75
$ws.Cells.Item(132, 13).Value = -755720
$ws.Cells.Item(132, 14).Value = -768557.75

$ws = $wb.Worksheets.Item("LTW")
# LTW row 7
$ws.Cells.Item(7, 8).Value = 3243.6
$ws.Cells.Item(7, 9).Value = 2805.2144
$ws.Cells.Item(7, 10).Value = 4266.5
$ws.Cells.Item(7, 11).Value = 2805.2144
$ws.Cells.Item(7, 12).Value = 4266.5
$ws.Cells.Item(7, 13).Value = -2693.2144
$ws.Cells.Item(7, 14).Value = -4490.5

# LTW row 40
$ws.Cells.Item(40, 8).Value = 5101.7715
$ws.Cells.Item(40, 9).Value = 5644.9473
$ws.Cells.Item(40, 10).Value = 4456.75
$ws.Cells.Item(40, 11).Value = 5644.9473
$ws.Cells.Item(40, 12).Value = 4456.75
$ws.Cells.Item(40, 13).Value = -5508.9473
$ws.Cells.Item(40, 14).Value = -4728.75

# LTW row 126
$ws.Cells.Item(126, 8).Value = 3243.6
$ws.Cells.Item(126, 9).Value = 2805.2144
$ws.Cells.Item(126, 10).Value = 4266.5
$ws.Cells.Item(126, 11).Value = 8415.643199999999
$ws.Cells.Item(126, 12).Value = 12799.5
$ws.Cells.Item(126, 13).Value = -5945.643199999999
$ws.Cells.Item(126, 14).Value = -17739.5

$ws = $wb.Worksheets.Item("WVR")
# WVR row 96
$ws.Cells.Item(96, 8).Value = 7918
$ws.Cells.Item(96, 9).Value = 4366.6665
$ws.Cells.Item(96, 11).Value = 4366.6665
$ws.Cells.Item(96, 13).Value = -2993.6665

# WVR row 109
$ws.Cells.Item(109, 8).Value = 24666.334
$ws.Cells.Item(109, 10).Value = 24666.334
$ws.Cells.Item(109, 12).Value = 24666.334
$ws.Cells.Item(109, 14).Value = -27440.334
